$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header translations: English transliterations -> uppercase Greek names
# (A1 "No" stays as-is)
$ws.Range("B1").Value = "ΜΑΚΑΡΙΟ"
$ws.Range("C1").Value = "ΝΕΟ ΓΣΠ"
$ws.Range("D1").Value = "ΑΓΛΑΝΤΖΙΑ"
$ws.Range("E1").Value = "ΛΑΤΣΙΑ"
$ws.Range("F1").Value = "ΕΡΓΑΤΕΣ"
$ws.Range("G1").Value = "ΑΝΑΓΥΙΑ"
$ws.Range("H1").Value = "ΝΗΣΟΥ"
$ws.Range("I1").Value = "ΔΑΛΙ"
$ws.Range("J1").Value = "ΛΥΜΠΙΑ"
$ws.Range("K1").Value = "ΚΟΚΚΙΝΟΤΡΙΜΙΘΙΑ"
$ws.Range("L1").Value = "ΠΑΛΙΟΜΕΤΟΧΟ"
$ws.Range("M1").Value = "ΜΑΜΜΑΡΙ"
$ws.Range("N1").Value = "ΠΕΡΙΣΤΕΡΩΝΑ"
$ws.Range("O1").Value = "ΛΑΚΑΤΑΜΕΙΑ"
$ws.Range("P1").Value = "ΚΟΤΣΙΑΤΗΣ"
$ws.Range("Q1").Value = "ΛΥΘΡΟΔΟΝΤΑΣ"
$ws.Range("R1").Value = "ΚΟΡΑΚΟΥ"
$ws.Range("S1").Value = "ΚΥΠΕΡΟΥΝΤΑ"
$ws.Range("T1").Value = "ΠΕΛΕΝΔΡΙ"
$ws.Range("U1").Value = "ΛΕΜΕΣΟΣ"
$ws.Range("V1").Value = "ΠΟΛΕΜΙΔΙΑ"
$ws.Range("W1").Value = "ΖΑΚΑΚΙ"
$ws.Range("X1").Value = "ΠΑΡΕΚΚΛΗΣΣΙΑ"
$ws.Range("Y1").Value = "ΚΟΛΟΣΣΙ"
$ws.Range("Z1").Value = "ΤΡΑΧΩΝΙ"
$ws.Range("AA1").Value = "ΕΡΗΜΗ"
$ws.Range("AB1").Value = "ΑΚΡΩΤΗΡΙ"
$ws.Range("AC1").Value = "ΓΕΡΜΑΣΟΓΕΙΑ"
$ws.Range("AD1").Value = "ΎΨΩΝΑΣ"
$ws.Range("AE1").Value = "ΛΑΡΝΑΚΑ"
$ws.Range("AF1").Value = "ΑΡΑΔΙΠΠΟΥ"
$ws.Range("AG1").Value = "ΑΓ. ΘΕΟΔΩΡΟΣ Λ/ΚΑΣ"
$ws.Range("AH1").Value = "ΨΕΥΔΑΣ"
$ws.Range("AI1").Value = "ΟΡΟΚΛΙΝΗ"
$ws.Range("AJ1").Value = "ΠΥΛΑ"
$ws.Range("AK1").Value = "ΚΟΡΝΟΣ"
$ws.Range("AL1").Value = "ΑΘΗΑΙΝΟΥ"
$ws.Range("AM1").Value = "ΛΕΙΒΑΔΙΑ"
$ws.Range("AN1").Value = "ΤΡΟΥΛΛΟΙ"
$ws.Range("AO1").Value = "ΚΙΤΙ"
$ws.Range("AP1").Value = "ΟΡΜΗΔΕΙΑ"
$ws.Range("AQ1").Value = "ΛΕΥΚΑΡΑ"
$ws.Range("AR1").Value = "ΞΥΛΟΦΑΓΟΥ"
$ws.Range("AS1").Value = "ΥΛΟΤΥΜΠΟΥ"
$ws.Range("AT1").Value = "ΔΑΣΑΚΙ ΑΧΝΑΣ"
$ws.Range("AU1").Value = "ΑΓΙΑ ΝΑΠΑ"
$ws.Range("AV1").Value = "ΔΕΡΥΝΕΙΑ"
$ws.Range("AW1").Value = "ΠΑΡΑΛΙΜΝΙ"
$ws.Range("AX1").Value = "ΑΥΓΟΡΟΥ"
$ws.Range("AY1").Value = "ΛΙΟΠΕΤΡΙ"
$ws.Range("AZ1").Value = "ΦΡΕΝΑΡΟΣ"
$ws.Range("BA1").Value = "ΣΩΤΗΡΑ"
$ws.Range("BB1").Value = "ΠΑΦΟΣ"
$ws.Range("BC1").Value = "ΧΛΩΡΑΚΑ"
$ws.Range("BD1").Value = "ΠΕΓΕΙΑ"
$ws.Range("BE1").Value = "ΓΕΡΟΣΚΗΠΟΥ"
$ws.Range("BF1").Value = "ΠΟΛΗΣ ΧΡΥΣΟΧΟΥΣ"
$ws.Range("BG1").Value = "ΑΓ. ΜΑΡΙΝΑ ΧΡΥΣΟΧΟΥΣ"

# Update the active selection to match the new state (header row selected)
$ws.Range("B1:BG1").Select() | Out-Null
